# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" on the rows that
#    reference it (Overview!E2:F2, zh-cn!C2, de-de!C2). Excel de-dupes
#    identical literal strings into a single shared-string entry, so writing
#    the same new text to all four cells keeps them collapsed onto one
#    shared string, exactly like the original file.
#
# 2. The "zh-cn"/"de-de" status columns on the Overview sheet (E:F) and the
#    "Status" column (C) on the per-locale sheets are narrowed from their old
#    auto-fit width down to the new auto-fit width produced by the shorter
#    "In Translation" header re-generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status values -------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the status columns -------------------------------------------------
# Excel's ColumnWidth is quantized to a whole-pixel grid (character units,
# 6pt steps) on write/read, same as interactive resizing in the real app;
# 12.5 is the input that lands on the pixel step closest to the new
# auto-fit width used by the report generator.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
